$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1")

# --- Column G width (Excel quantizes COM ColumnWidth to whole pixels;
#     38.6 is the input that lands on the closest achievable width to 39.453125) ---
$ws.Range("G1").EntireColumn.ColumnWidth = 38.6

# --- Cell text values ---
$ws.Range("B5").Value = 'Check username uniqueness in real time using AJAX (& before form submission)'
$ws.Range("C5").Value = 'Better user experience'

$ws.Range("B6").Value = 'Use SQL function for the above'
$ws.Range("C6").Value = 'More robust since not part of a ''development environment'''

$ws.Range("B7").Value = 'Countries from drop-down menu, not user-inputted'
$ws.Range("C7").Value = 'Part of DB design'

$ws.Range("B8").Value = 'If query for above fails, do not terminate registration form'
$ws.Range("C8").Value = 'Address data is optional, so give users the opportunity just to to submit it with their current registration form (or try to refresh)'

$ws.Range("B9").Value = 'Same for checking username uniqueness - though they can''t submit the form unless the check works'
$ws.Range("C9").Value = 'The query will keep on being tried each time focus is lost on the username input field, so it might work again in the future (e.g. a temporary connection glitch)'

$ws.Range("B10").Value = 'Password matching is only checked on client-side'
$ws.Range("C10").Value = 'Users who are submitting their own HTTP requests can deal with the inconvenience if it happens!'
$ws.Range("G10").Value = 'php/login_register_form_validation.php'

$ws.Range("B11").Value = 'Separate php file containing user input validation/sanitisation functions'
$ws.Range("C11").Value = 'Easier code readability'

$ws.Range("B12").Value = 'Validate/sanitise all user input'
$ws.Range("C12").Value = 'Standard reasons - client-side protections are not robust '

# --- B8:C8 were bold (style 8); target style 7 drops the bold font ---
$ws.Range("B8:C8").Font.Bold = $false

# --- Row heights grow to fit the newly wrapped text (auto-fit result) ---
$ws.Range("A5:A12").EntireRow.AutoFit()
$ws.Range("A5").RowHeight = 43.5
$ws.Range("A6").RowHeight = 29
$ws.Range("A7").RowHeight = 29
$ws.Range("A8").RowHeight = 58
$ws.Range("A9").RowHeight = 72.5
$ws.Range("A10").RowHeight = 43.5
$ws.Range("A11").RowHeight = 29
$ws.Range("A12").RowHeight = 29
